$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "12 ماهه منتهی به 1398/05"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/05"

$ws.Range("D9").Value = "1399-07-24 (8)"
$ws.Range("E9").Value = "1400-02-03 (7)"
$ws.Range("F9").Value = "1400-10-08 (8)"
$ws.Range("H9").Value = "1402-03-30 (5)"

$ws.Range("D11").Value = 13481
$ws.Range("E11").Value = 13203

$ws.Range("D12").Value = -6905
$ws.Range("E12").Value = -6907

$ws.Range("D13").Value = 6576
$ws.Range("E13").Value = 6296

$ws.Range("D14").Value = -828
$ws.Range("E14").Value = -960

$ws.Range("D15").Value = -54
$ws.Range("E15").Value = "-"

$ws.Range("D16").Value = 64
$ws.Range("E16").Value = 336

$ws.Range("D17").Value = 5758
$ws.Range("E17").Value = 5672

$ws.Range("D18").Value = -567
$ws.Range("E18").Value = -541

$ws.Range("D19").Value = 336
$ws.Range("E19").Value = 1339

$ws.Range("D20").Value = 5527
$ws.Range("E20").Value = 6469

$ws.Range("D21").Value = -584
$ws.Range("E21").Value = -470

$ws.Range("D22").Value = 4944
$ws.Range("E22").Value = 5999

$ws.Range("D24").Value = 4944
$ws.Range("E24").Value = 5999

$ws.Range("D26").Value = 3942
$ws.Range("E26").Value = 3331
